# Updated symbol list on Tue Dec 13 06:00:18 UTC 2022 with GitHub Actions
#
# The "Hora" column (G) moves from hour 5 to hour 6 for every data row,
# and several "Price" values (D) are refreshed with newer quotes.
# All of these cells are stored as text in the workbook, so values are
# assigned with a leading apostrophe to force text entry and avoid Excel
# reinterpreting numeric-looking strings as floating point numbers (which
# would corrupt exact decimal representations such as trailing zeros).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Price" (column D) values, keyed by row number. Only rows whose
# price actually changed are listed here; other rows keep their price.
$prices = @{
    2  = "267.72"
    3  = "21.40"
    4  = "6.265"
    5  = "0.06201"
    6  = "3.569"
    7  = "6.534"
    9  = "0.8236"
    10 = "0.1632"
    11 = "0.08210"
    12 = "0.03567"
    13 = "0.03153"
    14 = "0.09198"
    15 = "3.775"
    16 = "0.001629"
    17 = "0.04642"
    18 = "0.006424"
    19 = "0.006190"
    21 = "0.0001500"
    22 = "3.723"
    23 = "2.235"
    25 = "0.3320"
    28 = "0.0002713"
    40 = "0.04720"
    41 = "0.006926"
    42 = "0.004000"
    43 = "0.1121"
    45 = "0.00006318"
    46 = "0.0009901"
    48 = "0.9801"
    49 = "0.001141"
    50 = "0.00001900"
}

foreach ($row in $prices.Keys) {
    $ws.Range("D$row").Value = "'" + $prices[$row]
}

# "Hora" (column G) updates from 5 to 6 for every data row (2 through 51).
for ($row = 2; $row -le 51; $row++) {
    $ws.Range("G$row").Value = "'6"
}
